# Applies numeric corrections to several per-class Kujata_Profits worksheets
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR), matching the upstream commit diff.
$wb = $excel.ActiveWorkbook

# source diff hunk @@ -7022 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 835.3333
$ws.Range("I129").Value = 468.375
$ws.Range("J129").Value = 908.725
$ws.Range("K129").Value = 1405.125
$ws.Range("L129").Value = 2726.175
$ws.Range("M129").Value = 3594.875
$ws.Range("N129").Value = -12726.175

# source diff hunk @@ -7475 (ALC)
$ws.Range("H138").Value = 822029.7
$ws.Range("I138").Value = 1817.5
$ws.Range("J138").Value = 1012776.75
$ws.Range("K138").Value = 5452.5
$ws.Range("L138").Value = 3038330.25
$ws.Range("M138").Value = -312.5
$ws.Range("N138").Value = -3048610.25

# source diff hunk @@ -7915 (ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 192.83333
$ws.Range("I5").Value = 163
$ws.Range("J5").Value = 252.5
$ws.Range("K5").Value = 163
$ws.Range("L5").Value = 252.5
$ws.Range("M5").Value = -51
$ws.Range("N5").Value = -476.5

# source diff hunk @@ -8773 (ARM)
$ws.Range("H23").Value = 68754.75
$ws.Range("J23").Value = 57503.5
$ws.Range("L23").Value = 57503.5
$ws.Range("N23").Value = -58021.5

# source diff hunk @@ -9208 (ARM)
$ws.Range("H32").Value = 5399.0586
$ws.Range("I32").Value = 5611.5312
$ws.Range("J32").Value = 1999.5
$ws.Range("K32").Value = 5611.5312
$ws.Range("L32").Value = 1999.5
$ws.Range("M32").Value = -5324.5312
$ws.Range("N32").Value = -2573.5

# source diff hunk @@ -13582 (ARM)
$ws.Range("H122").Value = 2845.2
$ws.Range("I122").Value = 2570.6667
$ws.Range("J122").Value = 3257
$ws.Range("K122").Value = 7712.000100000001
$ws.Range("L122").Value = 9771
$ws.Range("M122").Value = -5262.000100000001
$ws.Range("N122").Value = -14671

# source diff hunk @@ -14063 (ARM)
$ws.Range("H132").Value = 2872.2
$ws.Range("I132").Value = 2511.476
$ws.Range("J132").Value = 3713.889
$ws.Range("K132").Value = 7534.428
$ws.Range("L132").Value = 11141.667
$ws.Range("M132").Value = -5004.428
$ws.Range("N132").Value = -16201.667

# source diff hunk @@ -14751 (BSM)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 192.83333
$ws.Range("I4").Value = 163
$ws.Range("J4").Value = 252.5
$ws.Range("K4").Value = 163
$ws.Range("L4").Value = 252.5
$ws.Range("M4").Value = -48
$ws.Range("N4").Value = -482.5

# source diff hunk @@ -15477 (BSM)
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

# source diff hunk @@ -20464 (BSM)
$ws.Range("H122").Value = 40000
$ws.Range("J122").Value = 40000
$ws.Range("L122").Value = 40000
$ws.Range("N122").Value = -49800

# source diff hunk @@ -22533 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 175105.25
$ws.Range("I22").Value = 140.33333
$ws.Range("J22").Value = 700000
$ws.Range("K22").Value = 140.33333
$ws.Range("L22").Value = 700000
$ws.Range("M22").Value = 209.66667
$ws.Range("N22").Value = -700700

# source diff hunk @@ -26052 (CRP)
$ws.Range("H94").Value = 890
$ws.Range("I94").Value = 670
$ws.Range("K94").Value = 670
$ws.Range("M94").Value = -219

# source diff hunk @@ -27415 (CRP)
$ws.Range("H122").Value = 842.2727
$ws.Range("I122").Value = 876.5
$ws.Range("J122").Value = 500
$ws.Range("K122").Value = 2629.5
$ws.Range("L122").Value = 1500
$ws.Range("M122").Value = -179.5
$ws.Range("N122").Value = -6400

# source diff hunk @@ -27908 (CRP)
$ws.Range("H132").Value = 2364.5
$ws.Range("I132").Value = 1710.5
$ws.Range("K132").Value = 5131.5
$ws.Range("M132").Value = -2601.5

# source diff hunk @@ -31213 (CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5801.4614
$ws.Range("I56").Value = 5801.4614
$ws.Range("K56").Value = 5801.4614
$ws.Range("M56").Value = -5271.4614

# source diff hunk @@ -35041 (CUL)
$ws.Range("H131").Value = 16396513
$ws.Range("I131").Value = 100000350
$ws.Range("J131").Value = 3603.255
$ws.Range("K131").Value = 300001050
$ws.Range("L131").Value = 10809.765
$ws.Range("M131").Value = -299996010
$ws.Range("N131").Value = -20889.765

# source diff hunk @@ -35457 (CUL)
$ws.Range("H139").Value = 1627.1702
$ws.Range("I139").Value = 1612.5217
$ws.Range("J139").Value = 1641.2084
$ws.Range("K139").Value = 4837.5651
$ws.Range("L139").Value = 4923.6252
$ws.Range("M139").Value = 302.4349000000002
$ws.Range("N139").Value = -15203.6252

# source diff hunk @@ -35704 (GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 274.81818
$ws.Range("I2").Value = 262.66666
$ws.Range("J2").Value = 289.4
$ws.Range("K2").Value = 262.66666
$ws.Range("L2").Value = 289.4
$ws.Range("M2").Value = -149.66666
$ws.Range("N2").Value = -515.4

# source diff hunk @@ -41533 (GSM)
$ws.Range("H122").Value = 8622566
$ws.Range("I122").Value = 1963.7368
$ws.Range("J122").Value = 25001708
$ws.Range("K122").Value = 5891.2104
$ws.Range("L122").Value = 75005124
$ws.Range("M122").Value = -3441.2104
$ws.Range("N122").Value = -75010024

# source diff hunk @@ -42020 (GSM)
$ws.Range("H132").Value = 2551.7334
$ws.Range("I132").Value = 2382.1904
$ws.Range("J132").Value = 2947.3333
$ws.Range("K132").Value = 7146.5712
$ws.Range("L132").Value = 8841.999899999999
$ws.Range("M132").Value = -4616.5712
$ws.Range("N132").Value = -13901.9999

# source diff hunk @@ -43575 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 812.4
$ws.Range("I22").Value = 544
$ws.Range("K22").Value = 544
$ws.Range("M22").Value = -249

# source diff hunk @@ -43820 (LTW)
$ws.Range("H27").Value = 812.4
$ws.Range("I27").Value = 544
$ws.Range("K27").Value = 544
$ws.Range("M27").Value = -437

# source diff hunk @@ -44454 (LTW)
$ws.Range("H40").Value = 2310.3333
$ws.Range("I40").Value = 2200.4443
$ws.Range("J40").Value = 2969.6667
$ws.Range("K40").Value = 2200.4443
$ws.Range("L40").Value = 2969.6667
$ws.Range("M40").Value = -2064.4443
$ws.Range("N40").Value = -3241.6667

# source diff hunk @@ -44745 (LTW)
$ws.Range("H46").Value = 3942.1538
$ws.Range("I46").Value = 725.8
$ws.Range("J46").Value = 5952.375
$ws.Range("K46").Value = 725.8
$ws.Range("L46").Value = 5952.375
$ws.Range("M46").Value = -537.8
$ws.Range("N46").Value = -6328.375

# source diff hunk @@ -47337 (LTW)
$ws.Range("H100").Value = 1264.8334
$ws.Range("I100").Value = 1247.25
$ws.Range("J100").Value = 1300
$ws.Range("K100").Value = 1247.25
$ws.Range("L100").Value = 1300
$ws.Range("M100").Value = -706.25
$ws.Range("N100").Value = -2382

# source diff hunk @@ -54011 (WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2397.1538
$ws.Range("I96").Value = 1999
$ws.Range("J96").Value = 2861.6667
$ws.Range("K96").Value = 1999
$ws.Range("L96").Value = 2861.6667
$ws.Range("M96").Value = -626
$ws.Range("N96").Value = -5607.6667

# source diff hunk @@ -54207 (WVR)
$ws.Range("H100").Value = 484.33334
$ws.Range("I100").Value = 475
$ws.Range("K100").Value = 950
$ws.Range("M100").Value = -409

# source diff hunk @@ -54841 (WVR)
$ws.Range("H113").Value = 301.79166
$ws.Range("I113").Value = 216.25
$ws.Range("K113").Value = 648.75
$ws.Range("M113").Value = 1521.25
